$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.667.87"
$ws.Range("E2").Value = "  -3.00%  "
$ws.Range("D3").Value = "2.086.28"
$ws.Range("E3").Value = "  -1.19%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.39%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "345.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.009"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5156"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.94%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4387"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09191"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "51.56"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.176"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.40"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.20%  "
$ws.Range("D13").Value = "2.091.88"
$ws.Range("E13").Value = "  -0.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.737"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.153"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001159"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.011"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "20.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06669"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.55%  "
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.190"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.04%  "
$ws.Range("D23").Value = "29.734.77"
$ws.Range("E23").Value = "  -2.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.303"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.39%  "
$ws.Range("D26").Value = "2.332.78"
$ws.Range("E26").Value = "  -1.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.523"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.144"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1052"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.635"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.200"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.941"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.134"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.16"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("E38").Value = "  -2.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06720"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2276"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.43"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6856"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  +1.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6644"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.300"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.622"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.52%  "
$ws.Range("E48").Value = "  -3.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000339"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "81.38"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07097"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.57%  "
